$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "r585"
$ws.Range("B10").Value = "tom"
$ws.Range("C10").Value = "Great news, turns out we not longer have high attenuation ever!"
$ws.Range("D10").Value = "2025-10-01 14:34:38"
